# "add README.md, fix database" -- the database-facing part of this commit
# replaces the old "pelanggan" (customer) sample sheet with a small
# "spare part" table: just NO PART / DESKRIPSI columns and a single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Wipe the old sample data (columns A:Q, rows 1:2) so stale shared strings
# don't linger and the sheet starts clean.
$ws.Range("A1:Q2").Value = $null

# New header row.
$ws.Range("A1").Value = "NO PART"
$ws.Range("B1").Value = "DESKRIPSI"

# New data row.
$ws.Range("A2").Value = 90111
$ws.Range("B2").Value = "gelo"
